$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1 (header row)
$ws.Range("F1").Value = "Last status check on: 13.01.2022 01:15"

# Update row 5 (Makro) price data
$ws.Range("B5").Value = 33.9
$ws.Range("C5").Value = 34.5

# D5 becomes a text cell holding "-0.6" (not a number)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "-0.6"
$ws.Range("D5").Style = "Normal"

# E5 becomes a text cell holding the timestamp string (date style removed)
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-01-13 01:15:09"
$ws.Range("E5").Style = "Normal"
